$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data for rows 2-7 (columns A-G)
$data = @(
    @(0, 2, "SMART SENSING MIDDLEWARE", 109.12, 1, 100, 3),
    @(1, 2, "SHAMIYANA APP", 79, 0.8, 63.2, 3),
    @(2, 2, "RAPID", 72.38, 1, 72.38, 3),
    @(3, 3, "Website for the Literature Society of the college", 98.59999999999999, 0.8, 78.88, 3),
    @(4, 3, "Post-processing of Large Language Models", 82.40000000000001, 1, 82.40000000000001, 3),
    @(5, 3, "Multi Model Data Analysis for Annotation of Human Activities", 69.16, 1, 69.16, 3)
)

$rowIndex = 2
foreach ($rec in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rec[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rec[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rec[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rec[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rec[4]
    $ws.Cells.Item($rowIndex, 6).Value = $rec[5]
    $ws.Cells.Item($rowIndex, 7).Value = $rec[6]
    $rowIndex++
}

# Remove the now-unused trailing rows (old rows 8 and 9)
$ws.Range("A8:G9").Delete()
